$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1: Camarilla pivot calculations (rows 3-14, columns C-K)
$ws1.Range("C3").Value = 4857.4
$ws1.Range("D3").Value = 197.44
$ws1.Range("E3").Value = 45269.56
$ws1.Range("F3").Value = 68188.32000000001
$ws1.Range("G3").Value = 689.97
$ws1.Range("H3").Value = 220.84
$ws1.Range("I3").Value = 161.98
$ws1.Range("J3").Value = 1184.44
$ws1.Range("K3").Value = 173.11

$ws1.Range("C4").Value = 4843.7
$ws1.Range("D4").Value = 195.87
$ws1.Range("E4").Value = 45151.13
$ws1.Range("F4").Value = 67880.81
$ws1.Range("G4").Value = 688.2
$ws1.Range("H4").Value = 219.78
$ws1.Range("I4").Value = 161.43
$ws1.Range("J4").Value = 1179.92
$ws1.Range("K4").Value = 172.74

$ws1.Range("C5").Value = 4830
$ws1.Range("D5").Value = 194.3
$ws1.Range("E5").Value = 45032.7
$ws1.Range("F5").Value = 67573.3
$ws1.Range("G5").Value = 686.4400000000001
$ws1.Range("H5").Value = 218.73
$ws1.Range("I5").Value = 160.89
$ws1.Range("J5").Value = 1175.4
$ws1.Range("K5").Value = 172.38

$ws1.Range("C6").Value = 4813.5
$ws1.Range("D6").Value = 192.4
$ws1.Range("E6").Value = 44891.35
$ws1.Range("F6").Value = 67208.64999999999
$ws1.Range("G6").Value = 684.3200000000001
$ws1.Range("H6").Value = 217.46
$ws1.Range("I6").Value = 160.25
$ws1.Range("J6").Value = 1169.9
$ws1.Range("K6").Value = 171.94

$ws1.Range("C7").Value = 4807.98
$ws1.Range("D7").Value = 191.76
$ws1.Range("E7").Value = 44844.06
$ws1.Range("F7").Value = 67086.66
$ws1.Range("G7").Value = 683.61
$ws1.Range("H7").Value = 217.04
$ws1.Range("I7").Value = 160.03
$ws1.Range("J7").Value = 1168.06
$ws1.Range("K7").Value = 171.79

$ws1.Range("C8").Value = 4802.5
$ws1.Range("D8").Value = 191.13
$ws1.Range("E8").Value = 44797.08
$ws1.Range("F8").Value = 66965.46000000001
$ws1.Range("G8").Value = 682.91
$ws1.Range("H8").Value = 216.62
$ws1.Range("I8").Value = 159.82
$ws1.Range("J8").Value = 1166.23
$ws1.Range("K8").Value = 171.65

$ws1.Range("C9").Value = 4791.5
$ws1.Range("D9").Value = 189.87
$ws1.Range("E9").Value = 44702.92
$ws1.Range("F9").Value = 66722.53999999999
$ws1.Range("G9").Value = 681.49
$ws1.Range("H9").Value = 215.78
$ws1.Range("I9").Value = 159.38
$ws1.Range("J9").Value = 1162.57
$ws1.Range("K9").Value = 171.35

$ws1.Range("C10").Value = 4786.02
$ws1.Range("D10").Value = 189.24
$ws1.Range("E10").Value = 44655.94
$ws1.Range("F10").Value = 66601.34
$ws1.Range("G10").Value = 680.79
$ws1.Range("H10").Value = 215.36
$ws1.Range("I10").Value = 159.17
$ws1.Range("J10").Value = 1160.74
$ws1.Range("K10").Value = 171.21

$ws1.Range("C11").Value = 4780.5
$ws1.Range("D11").Value = 188.6
$ws1.Range("E11").Value = 44608.65
$ws1.Range("F11").Value = 66479.35000000001
$ws1.Range("G11").Value = 680.08
$ws1.Range("H11").Value = 214.94
$ws1.Range("I11").Value = 158.95
$ws1.Range("J11").Value = 1158.9
$ws1.Range("K11").Value = 171.06

$ws1.Range("C12").Value = 4764
$ws1.Range("D12").Value = 186.7
$ws1.Range("E12").Value = 44467.3
$ws1.Range("F12").Value = 66114.7
$ws1.Range("G12").Value = 677.96
$ws1.Range("H12").Value = 213.67
$ws1.Range("I12").Value = 158.31
$ws1.Range("J12").Value = 1153.4
$ws1.Range("K12").Value = 170.62

$ws1.Range("C13").Value = 4750.3
$ws1.Range("D13").Value = 185.13
$ws1.Range("E13").Value = 44348.87
$ws1.Range("F13").Value = 65807.19
$ws1.Range("G13").Value = 676.2
$ws1.Range("H13").Value = 212.62
$ws1.Range("I13").Value = 157.77
$ws1.Range("J13").Value = 1148.88
$ws1.Range("K13").Value = 170.26

$ws1.Range("C14").Value = 4736.6
$ws1.Range("D14").Value = 183.56
$ws1.Range("E14").Value = 44230.44
$ws1.Range("F14").Value = 65499.68
$ws1.Range("G14").Value = 674.4299999999999
$ws1.Range("H14").Value = 211.56
$ws1.Range("I14").Value = 157.22
$ws1.Range("J14").Value = 1144.36
$ws1.Range("K14").Value = 169.89

# Sheet2: raw Open/High/Low/Close/Previous Close data (rows 2-10, columns C-G)
$ws2.Range("C2").Value = 4778
$ws2.Range("D2").Value = 4825
$ws2.Range("E2").Value = 4765
$ws2.Range("F2").Value = 4797
$ws2.Range("G2").Value = 4796

$ws2.Range("C3").Value = 194.3
$ws2.Range("D3").Value = 196.2
$ws2.Range("E3").Value = 189.3
$ws2.Range("F3").Value = 190.5
$ws2.Range("G3").Value = 193.8

$ws2.Range("C4").Value = 44741
$ws2.Range("D4").Value = 44785
$ws2.Range("E4").Value = 44271
$ws2.Range("F4").Value = 44750
$ws2.Range("G4").Value = 44879

$ws2.Range("C5").Value = 67200
$ws2.Range("D5").Value = 67259
$ws2.Range("E5").Value = 65933
$ws2.Range("F5").Value = 66844
$ws2.Range("G5").Value = 67545

$ws2.Range("C6").Value = 683.8
$ws2.Range("D6").Value = 683.8
$ws2.Range("E6").Value = 676.1
$ws2.Range("F6").Value = 682.2
$ws2.Range("G6").Value = 686.55

$ws2.Range("C7").Value = 218.65
$ws2.Range("D7").Value = 219.1
$ws2.Range("E7").Value = 214.5
$ws2.Range("F7").Value = 216.2
$ws2.Range("G7").Value = 218.55

$ws2.Range("C8").Value = 158.7
$ws2.Range("D8").Value = 160.25
$ws2.Range("E8").Value = 157.9
$ws2.Range("F8").Value = 159.6
$ws2.Range("G8").Value = 158.8

$ws2.Range("C9").Value = 1180.3
$ws2.Range("D9").Value = 1182
$ws2.Range("E9").Value = 1162
$ws2.Range("F9").Value = 1164.4
$ws2.Range("G9").Value = 1186.2

$ws2.Range("C10").Value = 172.45
$ws2.Range("D10").Value = 172.45
$ws2.Range("E10").Value = 170.85
$ws2.Range("F10").Value = 171.5
$ws2.Range("G10").Value = 172.95
